# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 448, pushing the existing rows
# (448:468) down to (449:469). The new row holds a fresh "Ajo" / "Chino"
# price entry for Femacal de La Calera - Coquimbo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 448 - this shifts rows
# 448:468 down to 449:469, matching the diff exactly (old row 468's
# data ends up as the new row 469).
$ws.Rows("448:448").Insert()

# Populate the newly inserted row 448 with the new weekly record.
$ws.Range("A448").Value = 3
$ws.Range("B448").Value = "Femacal de La Calera"
$ws.Range("C448").Value = "Coquimbo"
$ws.Range("D448").Value = 44753
$ws.Range("E448").Value = 5
$ws.Range("F448").Value = 100112003
$ws.Range("G448").Value = "Ajo"
$ws.Range("H448").Value = "Chino"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 85
$ws.Range("K448").Value = 19000
$ws.Range("L448").Value = 20000
$ws.Range("M448").Value = 19471
$ws.Range("N448").Value = "$/caja 10 kilos"
$ws.Range("O448").Value = "China"
$ws.Range("P448").Value = 1947
$ws.Range("Q448").Value = 10
$ws.Range("R448").Value = "Hortaliza"
